$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Cells.Item(2,4).Value = "ECs"                    # D2
$ws.Cells.Item(2,5).Value = 3                          # E2
$ws.Cells.Item(2,7).Value = 4.062756666666667          # G2
$ws.Cells.Item(2,8).Value = 12.18827                   # H2
$ws.Cells.Item(2,9).Value = 0.6829811567947219         # I2
$ws.Cells.Item(2,10).Value = 0.6829811567947218        # J2
$ws.Cells.Item(2,11).Value = 2                         # K2
$ws.Cells.Item(2,12).Value = 0.6666666666666666        # L2
$ws.Cells.Item(2,13).Value = 83.95844533333333         # M2
$ws.Cells.Item(2,14).Value = 251.875336                # N2
$ws.Cells.Item(2,15).Value = 0.9979754487867319        # O2
$ws.Cells.Item(2,16).Value = 0.9979754487867319        # P2
$ws.Cells.Item(2,17).Value = 341.1027335009689         # Q2
$ws.Cells.Item(2,18).Value = 3069.92460150872          # R2
$ws.Cells.Item(2,19).Value = 0.6815984264650939        # S2
$ws.Cells.Item(2,20).Value = 0.6815984264650938        # T2

# --- Row 3 updates ---
$ws.Cells.Item(3,1).Value = "ECs"                      # A3
$ws.Cells.Item(3,5).Value = 3                          # E3
$ws.Cells.Item(3,7).Value = 4.062756666666667          # G3
$ws.Cells.Item(3,8).Value = 12.18827                   # H3
$ws.Cells.Item(3,9).Value = 0.6829811567947219         # I3
$ws.Cells.Item(3,10).Value = 0.6829811567947218        # J3
$ws.Cells.Item(3,11).Value = 3                         # K3
$ws.Cells.Item(3,13).Value = 0.170323                  # M3
$ws.Cells.Item(3,14).Value = 0.510969                  # N3
$ws.Cells.Item(3,15).Value = 0.002024551213268089      # O3
$ws.Cells.Item(3,16).Value = 0.00202455121326809       # P3
$ws.Cells.Item(3,17).Value = 0.6919809037366668        # Q3
$ws.Cells.Item(3,18).Value = 6.227828133630001         # R3
$ws.Cells.Item(3,19).Value = 0.001382730329627997      # S3
$ws.Cells.Item(3,20).Value = 0.001382730329627997      # T3

# --- Row 4 (new) ---
$ws.Cells.Item(4,1).Value = "FAPs"                     # A4
$ws.Cells.Item(4,2).Value = "Icam5"                    # B4
$ws.Cells.Item(4,3).Value = "Itgb2"                    # C4
$ws.Cells.Item(4,4).Value = "ECs"                      # D4
$ws.Cells.Item(4,5).Value = 3                          # E4
$ws.Cells.Item(4,6).Value = 1                          # F4
$ws.Cells.Item(4,7).Value = 1.885806666666667          # G4
$ws.Cells.Item(4,8).Value = 5.65742                    # H4
$ws.Cells.Item(4,9).Value = 0.3170188432052781         # I4
$ws.Cells.Item(4,10).Value = 0.3170188432052781        # J4
$ws.Cells.Item(4,11).Value = 2                         # K4
$ws.Cells.Item(4,12).Value = 0.6666666666666666        # L4
$ws.Cells.Item(4,13).Value = 83.95844533333333         # M4
$ws.Cells.Item(4,14).Value = 251.875336                # N4
$ws.Cells.Item(4,15).Value = 0.9979754487867319        # O4
$ws.Cells.Item(4,16).Value = 0.9979754487867319        # P4
$ws.Cells.Item(4,17).Value = 158.3293959325689         # Q4
$ws.Cells.Item(4,18).Value = 1424.96456339312          # R4
$ws.Cells.Item(4,19).Value = 0.3163770223216381        # S4
$ws.Cells.Item(4,20).Value = 0.3163770223216381        # T4

# --- Row 5 (new) ---
$ws.Cells.Item(5,1).Value = "FAPs"                     # A5
$ws.Cells.Item(5,2).Value = "Icam5"                    # B5
$ws.Cells.Item(5,3).Value = "Itgb2"                    # C5
$ws.Cells.Item(5,4).Value = "FAPs"                     # D5
$ws.Cells.Item(5,5).Value = 3                          # E5
$ws.Cells.Item(5,6).Value = 1                          # F5
$ws.Cells.Item(5,7).Value = 1.885806666666667          # G5
$ws.Cells.Item(5,8).Value = 5.65742                    # H5
$ws.Cells.Item(5,9).Value = 0.3170188432052781         # I5
$ws.Cells.Item(5,10).Value = 0.3170188432052781        # J5
$ws.Cells.Item(5,11).Value = 3                         # K5
$ws.Cells.Item(5,12).Value = 1                         # L5
$ws.Cells.Item(5,13).Value = 0.170323                  # M5
$ws.Cells.Item(5,14).Value = 0.510969                  # N5
$ws.Cells.Item(5,15).Value = 0.002024551213268089      # O5
$ws.Cells.Item(5,16).Value = 0.00202455121326809       # P5
$ws.Cells.Item(5,17).Value = 0.3211962488866666        # Q5
$ws.Cells.Item(5,18).Value = 2.89076623998             # R5
$ws.Cells.Item(5,19).Value = 0.000641820883640092      # S5
$ws.Cells.Item(5,20).Value = 0.0006418208836400921     # T5

Write-Host "edits applied"
